# Refresh the "Price" column (D) with updated cryptocurrency quotes.
#
# The source data is stored as literal text (not numbers) so that values
# such as "274.20", "3.660" or "0.0001500" keep their exact printed form
# (trailing zeros, fixed decimal widths, etc.). Assigning a numeric-looking
# string straight to Range.Value would make Excel coerce it into a real
# number (dropping the formatting-significant trailing zeros) and would
# also stamp the cell with a "number stored as text" style the moment we
# try to force it back to text. To avoid both side effects we:
#   1. assign the value with a leading apostrophe, which is exactly how a
#      user forces Excel to keep a numeric-looking entry as text; then
#   2. call ClearFormats() on the cell so the only trace left behind is
#      the text itself, not the "stored as text" quote-prefix styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = [ordered]@{
    "D2"  = "273.69"
    "D3"  = "22.93"
    "D4"  = "6.357"
    "D5"  = "0.06232"
    "D6"  = "3.660"
    "D7"  = "6.663"
    "D8"  = "1.374"
    "D9"  = "0.8316"
    "D10" = "0.01375"
    "D11" = "0.1631"
    "D12" = "0.08315"
    "D13" = "0.03400"
    "D14" = "0.03105"
    "D15" = "0.09309"
    "D16" = "3.893"
    "D17" = "0.001640"
    "D19" = "0.006343"
    "D20" = "0.005553"
    "D22" = "0.0001500"
    "D23" = "3.722"
    "D24" = "2.322"
    "D40" = "0.04685"
    "D41" = "0.007040"
    "D42" = "0.1164"
    "D43" = "0.003351"
    "D44" = "0.01215"
    "D45" = "0.00006260"
    "D47" = "0.9005"
    "D48" = "0.02927"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    $cell.ClearFormats()
}
